$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to numbers by Excel
# need NumberFormat forced to Text ("@") before assignment so they stay as strings,
# matching the inline-string cell type used throughout the sheet.

$ws.Range("D2").Value = '29.687.35'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '1.602.89'
$ws.Range("E3").Value = '  +1.36%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.36'
$ws.Range("E5").Value = '  -0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  +0.94%  '

$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '27.85'
$ws.Range("E8").Value = '  +5.45%  '

$ws.Range("E9").Value = '  +1.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0603'
$ws.Range("E10").Value = '  +1.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("E11").Value = '  +0.60%  '

$ws.Range("D12").Value = '1.832.26'
$ws.Range("E12").Value = '  +1.37%  '

$ws.Range("D13").Value = '1.597.01'
$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.546'
$ws.Range("E14").Value = '  +4.14%  '

$ws.Range("D15").Value = '29.678.81'
$ws.Range("E15").Value = '  +1.52%  '

$ws.Range("E16").Value = '  +0.81%  '

$ws.Range("E17").Value = '  +1.89%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.83'
$ws.Range("E18").Value = '  +1.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.80'
$ws.Range("E19").Value = '  +4.47%  '

$ws.Range("E20").Value = '  +1.20%  '

$ws.Range("E21").Value = '  +0.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.02'
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.45'
$ws.Range("E23").Value = '  +2.85%  '

$ws.Range("E24").Value = '  -0.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.54'
$ws.Range("E25").Value = '  +0.84%  '

$ws.Range("E26").Value = '  +1.98%  '

$ws.Range("E27").Value = '  +0.40%  '

$ws.Range("E28").Value = '  +0.94%  '

$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("E30").Value = '  +2.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.07'
$ws.Range("E31").Value = '  -0.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("E32").Value = '  +0.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.19'
$ws.Range("E33").Value = '  +3.59%  '

$ws.Range("D34").Value = '1.425.96'
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("B35").Value = 'MXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.92'
$ws.Range("E35").Value = '  +6.03%  '

$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.56'
$ws.Range("E36").Value = '  +3.53%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.03'
$ws.Range("E37").Value = '  -0.92%  '

$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("E39").Value = '  +2.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '58.15'
$ws.Range("E40").Value = '  +6.67%  '

$ws.Range("E41").Value = '  +2.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0500'
$ws.Range("E42").Value = '  +6.38%  '

$ws.Range("E43").Value = '  +0.33%  '

$ws.Range("E44").Value = '  +2.96%  '

$ws.Range("E45").Value = '  +0.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.44'
$ws.Range("E46").Value = '  +2.79%  '

$ws.Range("E47").Value = '  +16.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.36'
$ws.Range("E48").Value = '  +0.16%  '

$ws.Range("D49").Value = '1.743.15'
$ws.Range("E49").Value = '  +1.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.82'
$ws.Range("E50").Value = '  +1.56%  '

$ws.Range("E51").Value = '  +3.43%  '
